$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 69  # column BQ ("cMax")

function Get-CellValue($cell) {
    return $cell.Value()
}

function Set-CellValueSafe($cell, $val) {
    # A leading apostrophe is used by Excel as a "treat next as text" marker and is
    # swallowed when assigned through .Value. Double it so the literal leading
    # apostrophe in the data is preserved.
    if ($val -is [string] -and $val.Length -gt 0 -and $val.Substring(0,1) -eq "'") {
        $cell.Value = "'" + $val
    } else {
        $cell.Value = $val
    }
}

for ($c = 1; $c -le $lastCol; $c++) {
    $cell14 = $ws.Cells.Item(14, $c)
    $cell15 = $ws.Cells.Item(15, $c)

    $v14 = Get-CellValue $cell14
    $v15 = Get-CellValue $cell15

    Set-CellValueSafe $cell14 $v15
    Set-CellValueSafe $cell15 $v14
}
